$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AUW combine re-run: append the latest executed report row (row 12),
# copying the formatting already used by the existing data rows.
$ws.Range("A12").Value = 44573
$ws.Range("A12").NumberFormat = $ws.Range("A10").NumberFormat

$ws.Range("B12").Value = "Production"

$ws.Range("C12").Value = 151
$ws.Range("D12").Value = 149
$ws.Range("E12").Value = 2

$ws.Range("F12").Value = "After execution all test cases pass"
$ws.Range("G12").Value = "Test cases initially fail because of page load affected by network"

$ws.Rows.Item(12).RowHeight = 75

# Excel leaves a fresh blank row below the newly appended data (row 13),
# carrying the date column's number format down into it.
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat

# Leave the selection/scroll position the way Excel does after typing the
# new row and landing on the next blank one.
$ws.Range("A13:L13").Select()
$excel.ActiveWindow.ScrollRow = 10
